$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 22 (shifts old rows 22-24 down to 23-25) to make room
# for splitting the combined "track / attack" note into two separate rows.
$ws.Rows("22:22").Insert()

# ---------------------------------------------------------------------
# Fix up the Status-column (column C) cell formatting *before* touching
# the values, while the original fills are still sitting on known cells:
#   - C25 (shifted down from the old C24) still carries the "/" fill
#     -> copy that format onto the new C22 row.
#   - C21 still carries the "BUSY" fill
#     -> copy that format onto C25 (which needs to become BUSY).
#   - C2 carries the "DONE" fill, used for every DONE status cell.
# ---------------------------------------------------------------------
$ws.Range("C25").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("C21").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("C2").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 19 (unchanged content) ---
$ws.Range("A19").Value2 = "Let characters follow selected character ( navmesh )"
$ws.Range("B19").Value2 = 0.041666666666666664
$ws.Range("C19").Value2 = "DONE"

# --- Row 20 (unchanged content) ---
$ws.Range("A20").Value2 = "Basic enemy movement ( navmesh agents )"
$ws.Range("B20").Value2 = 0.006944444444444444
$ws.Range("C20").Value2 = "DONE"

# --- Row 21: first half of the split note, status now DONE ---
$ws.Range("A21").Value2 = "Let enemies track ( follow ) selected character ( navmesh ). "
$ws.Range("B21").Value2 = 0.08333333333333333
$ws.Range("C21").Value2 = "DONE"
$ws.Range("D21").Value2 = "Now enemy tracks always selected character, I need a check to see if one of the enemies is close enough"

# --- Row 22 (new row): second half of the split note, status "/" ---
$ws.Range("A22").Value2 = "When close enough, attack closest of 3 characters"
$ws.Range("C22").Value2 = "/"

# --- Row 23 (was row 22): unchanged content, just shifted down ---
$ws.Range("A23").Value2 = "Revisit end game, more elegant solution"
$ws.Range("B23").Value2 = 0.0625
$ws.Range("C23").Value2 = "DONE"
$ws.Range("D23").Value2 = "This cost me a lot of time due to testing"

# --- Row 24 (was row 23): unchanged content, just shifted down ---
$ws.Range("A24").Value2 = "Inventory ( panel )"
$ws.Range("B24").Value2 = 0.006944444444444444
$ws.Range("C24").Value2 = "DONE"
$ws.Range("D24").Value2 = "Seperated character panel and inventory panel ( put actual items on the inventory panel )"

# --- Row 25 (was row 24): status changes from "/" to "BUSY" ---
$ws.Range("A25").Value2 = "Cleanup code"
$ws.Range("C25").Value2 = "BUSY"

# Update the selection to match the final cursor position recorded in the file.
$ws.Range("B27").Select()
